$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New training-day column AL, date 2025-09-02 (serial 45902). Pick up the
# exact formatting (date number format, centered + vertically centered
# alignment) already used for the neighbouring header cell by copying it
# rather than re-declaring a style from scratch.
$ws.Range("AL1").Value2 = 45902
$ws.Range("M1").Copy()
$ws.Range("AL1").PasteSpecial(-4122)

# Per-player attendance marks for the new training day, mirroring the
# formatting already used on the rest of the row (column AK).
$marks = @{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "P"
    6  = "P"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "B"
    11 = "P"
    12 = "A"
    13 = "P"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "RH"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "P"
    22 = "P"
    23 = "P"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "RH"
}

foreach ($row in $marks.Keys) {
    $ws.Cells.Item($row, 37).Copy()
    $ws.Cells.Item($row, 38).PasteSpecial(-4122)
    $ws.Cells.Item($row, 38).Value = $marks[$row]
}

# Totals row: count of "P" marks for the new column.
$ws.Range("AL28").Formula = "=COUNTIF(AL2:AL27,""P"")"
